# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-12-01 08:24:01
#
# This script updates the "Recorded By" (G), "Students" (H), and the
# "Average Attendance %" (L10 / S15) cells on the "Session Analysis Results"
# sheet to reflect the refreshed attendance data pulled from the main
# repository.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell as literal text, without letting Excel's
# automatic value-parsing reinterpret number/fraction/percentage-looking
# strings (e.g. "102/251" or "26.0%") as dates, fractions, or percentages.
# Writing it as a quoted text formula forces Excel to store it as a string
# result; converting that formula to a static value (copy / paste-special
# values) afterwards drops the formula while keeping the literal text and
# the cell's original style untouched.
function Set-LiteralText {
    param(
        [string]$CellAddress,
        [string]$Text
    )

    $target = $ws.Range($CellAddress)
    $target.Formula = '="' + $Text + '"'
    $target.Copy()
    $target.PasteSpecial(-4163)  # xlPasteValues
}

# Row 2 - GIT session 1: reorder "Recorded By" list
$ws.Range("G2").Value = "gehanadel@med.asu.edu.eg, System, Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"

# Row 3 - GIT session 2: reorder "Recorded By" list
$ws.Range("G3").Value = "asmaa.reda@med.asu.edu.eg, System, majorelle.magdy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"

# Row 4 - GIT session 3: reorder "Recorded By" list
$ws.Range("G4").Value = "gehanadel@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"

# Row 5 - GIT session 4: updated "Recorded By" list (new recorder added) and refreshed attendance count
$ws.Range("G5").Value = "asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
Set-LiteralText "H5" "102/251"

# Row 6 - GIT session 5: reorder "Recorded By" list
$ws.Range("G6").Value = "majorelle.magdy@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, Mohammedeltanany@med.asu.edu.eg"

# Row 7 - BIOCHEMISTRY LAB/CBL: reorder "Recorded By" list
$ws.Range("G7").Value = "Amera.a.saad@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Fatmaelhady@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg"

# Row 10 - GIT group statistics: refreshed average attendance %
Set-LiteralText "L10" "26.0%"

# Row 12 - MICROBIOLOGY: reorder "Recorded By" list
$ws.Range("G12").Value = "Eman.m.abosakaya@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, dina.adel@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg"

# Row 15 - Liver group statistics: refreshed average attendance %
Set-LiteralText "S15" "26.0%"

# Row 27 - PHARMACOLOGY: reorder "Recorded By" list
$ws.Range("G27").Value = "nourhan.mostafa@med.asu.edu.eg, hana.amr@med.asu.edu.eg"
